# Daily attendance processing - 2025-10-09 23:18:53
#
# Normalizes the "Recorded By" column (G): whenever the recorder list
# begins with the literal "System" entry, that entry is moved from the
# front of the comma-separated list to the end, e.g.
#   "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#   "System, system, backup@backdoor.com" -> "system, backup@backdoor.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Range("G" + $r)
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ", "

        if ($parts[0] -eq "System") {
            $rest = $parts[1..($parts.Length - 1)]
            $newValue = ($rest + "System") -join ", "
            $cell.Value = $newValue
        }
    }
}
